$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5295576666666667
$ws.Range("H2").Value = 1.588673
$ws.Range("I2").Value = 0.7656712979474436
$ws.Range("J2").Value = 0.7656712979474436
$ws.Range("M2").Value = 1.979087666666667
$ws.Range("N2").Value = 5.937263000000001
$ws.Range("O2").Value = 0.2426859505365239
$ws.Range("P2").Value = 0.2426859505365239
$ws.Range("Q2").Value = 1.048041046888778
$ws.Range("R2").Value = 9.432369421999001
$ws.Range("S2").Value = 0.1858176667409094
$ws.Range("T2").Value = 0.1858176667409094

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5295576666666667
$ws.Range("H3").Value = 1.588673
$ws.Range("I3").Value = 0.7656712979474436
$ws.Range("J3").Value = 0.7656712979474436
$ws.Range("N3").Value = 7.706687
$ws.Range("O3").Value = 0.3150112535157145
$ws.Range("P3").Value = 0.3150112535157146
$ws.Range("Q3").Value = 1.360378395150111
$ws.Range("R3").Value = 12.243405556351
$ws.Range("S3").Value = 0.2411950753474283
$ws.Range("T3").Value = 0.2411950753474284

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5295576666666667
$ws.Range("H4").Value = 1.588673
$ws.Range("I4").Value = 0.7656712979474436
$ws.Range("J4").Value = 0.7656712979474436
$ws.Range("M4").Value = 2.580162
$ws.Range("N4").Value = 7.740486
$ws.Range("O4").Value = 0.3163927894931816
$ws.Range("P4").Value = 0.3163927894931816
$ws.Range("Q4").Value = 1.366344568342
$ws.Range("R4").Value = 12.297101115078
$ws.Range("S4").Value = 0.2422528777924566
$ws.Range("T4").Value = 0.2422528777924567

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5295576666666667
$ws.Range("H5").Value = 1.588673
$ws.Range("I5").Value = 0.7656712979474436
$ws.Range("J5").Value = 0.7656712979474436
$ws.Range("M5").Value = 1.026787666666667
$ws.Range("N5").Value = 3.080363
$ws.Range("O5").Value = 0.1259100064545799
$ws.Range("P5").Value = 0.1259100064545799
$ws.Range("Q5").Value = 0.5437432809221111
$ws.Range("R5").Value = 4.893689528299
$ws.Range("S5").Value = 0.09640567806664918
$ws.Range("T5").Value = 0.09640567806664921

# Row 6
$ws.Range("I6").Value = 0.1255607564018283
$ws.Range("J6").Value = 0.1255607564018283
$ws.Range("M6").Value = 1.979087666666667
$ws.Range("N6").Value = 5.937263000000001
$ws.Range("O6").Value = 0.2426859505365239
$ws.Range("P6").Value = 0.2426859505365239
$ws.Range("Q6").Value = 0.171865952061
$ws.Range("R6").Value = 1.546793568549
$ws.Range("S6").Value = 0.03047183151746264
$ws.Range("T6").Value = 0.03047183151746265

# Row 7
$ws.Range("I7").Value = 0.1255607564018283
$ws.Range("J7").Value = 0.1255607564018283
$ws.Range("N7").Value = 7.706687
$ws.Range("O7").Value = 0.3150112535157145
$ws.Range("P7").Value = 0.3150112535157146
$ws.Range("S7").Value = 0.03955305126652122
$ws.Range("T7").Value = 0.03955305126652123

# Row 8
$ws.Range("I8").Value = 0.1255607564018283
$ws.Range("J8").Value = 0.1255607564018283
$ws.Range("M8").Value = 2.580162
$ws.Range("N8").Value = 7.740486
$ws.Range("O8").Value = 0.3163927894931816
$ws.Range("P8").Value = 0.3163927894931816
$ws.Range("Q8").Value = 0.224063848242
$ws.Range("R8").Value = 2.016574634178
$ws.Range("S8").Value = 0.03972651796884832
$ws.Range("T8").Value = 0.03972651796884833

# Row 9
$ws.Range("I9").Value = 0.1255607564018283
$ws.Range("J9").Value = 0.1255607564018283
$ws.Range("M9").Value = 1.026787666666667
$ws.Range("N9").Value = 3.080363
$ws.Range("O9").Value = 0.1259100064545799
$ws.Range("P9").Value = 0.1259100064545799
$ws.Range("Q9").Value = 0.08916726776099999
$ws.Range("R9").Value = 0.8025054098490001
$ws.Range("S9").Value = 0.01580935564899614
$ws.Range("T9").Value = 0.01580935564899614

# Row 10
$ws.Range("G10").Value = 0.07522666666666666
$ws.Range("I10").Value = 0.108767945650728
$ws.Range("J10").Value = 0.108767945650728
$ws.Range("M10").Value = 1.979087666666667
$ws.Range("N10").Value = 5.937263000000001
$ws.Range("O10").Value = 0.2426859505365239
$ws.Range("P10").Value = 0.2426859505365239
$ws.Range("Q10").Value = 0.1488801682044444
$ws.Range("R10").Value = 1.33992151384
$ws.Range("S10").Value = 0.0263964522781519
$ws.Range("T10").Value = 0.02639645227815191

# Row 11
$ws.Range("G11").Value = 0.07522666666666666
$ws.Range("I11").Value = 0.108767945650728
$ws.Range("J11").Value = 0.108767945650728
$ws.Range("N11").Value = 7.706687
$ws.Range("O11").Value = 0.3150112535157145
$ws.Range("P11").Value = 0.3150112535157146
$ws.Range("S11").Value = 0.03426312690176495
$ws.Range("T11").Value = 0.03426312690176495

# Row 12
$ws.Range("G12").Value = 0.07522666666666666
$ws.Range("I12").Value = 0.108767945650728
$ws.Range("J12").Value = 0.108767945650728
$ws.Range("M12").Value = 2.580162
$ws.Range("N12").Value = 7.740486
$ws.Range("O12").Value = 0.3163927894931816
$ws.Range("P12").Value = 0.3163927894931816
$ws.Range("Q12").Value = 0.19409698672
$ws.Range("R12").Value = 1.74687288048
$ws.Range("S12").Value = 0.03441339373187661
$ws.Range("T12").Value = 0.03441339373187662

# Row 13
$ws.Range("G13").Value = 0.07522666666666666
$ws.Range("I13").Value = 0.108767945650728
$ws.Range("J13").Value = 0.108767945650728
$ws.Range("M13").Value = 1.026787666666667
$ws.Range("N13").Value = 3.080363
$ws.Range("O13").Value = 0.1259100064545799
$ws.Range("P13").Value = 0.1259100064545799
$ws.Range("Q13").Value = 0.07724181353777777
$ws.Range("R13").Value = 0.69517632184
$ws.Range("S13").Value = 0.01369497273893456
$ws.Range("T13").Value = 0.01369497273893456
